# Sync automatico del tracker: append nuevas filas de resultados de partidos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fila 75 - event_id 14357927
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = '14357927'
$ws.Range("A75").ClearFormats()
$ws.Range("B75").NumberFormat = "@"
$ws.Range("B75").Value = '2025-08-06'
$ws.Range("B75").ClearFormats()
$ws.Range("C75").Value = 'Arthur Fery'
$ws.Range("D75").Value = 'Martin Landaluce'
$ws.Range("E75").Value = 'Gana Martin Landaluce'
$ws.Range("F75").Value = 1.62
$ws.Range("G75").Value = "'"
$ws.Range("G75").ClearFormats()
$ws.Range("H75").Value = "'"
$ws.Range("H75").ClearFormats()

# Fila 76 - event_id 14350781
$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = '14350781'
$ws.Range("A76").ClearFormats()
$ws.Range("B76").NumberFormat = "@"
$ws.Range("B76").Value = '2025-08-06'
$ws.Range("B76").ClearFormats()
$ws.Range("C76").Value = 'Mikhail Kukushkin'
$ws.Range("D76").Value = 'Emilio Nava'
$ws.Range("E76").Value = 'Gana Mikhail Kukushkin'
$ws.Range("F76").Value = 4
$ws.Range("G76").Value = "'"
$ws.Range("G76").ClearFormats()
$ws.Range("H76").Value = "'"
$ws.Range("H76").ClearFormats()

# Fila 77 - event_id 14359046
$ws.Range("A77").NumberFormat = "@"
$ws.Range("A77").Value = '14359046'
$ws.Range("A77").ClearFormats()
$ws.Range("B77").NumberFormat = "@"
$ws.Range("B77").Value = '2025-08-06'
$ws.Range("B77").ClearFormats()
$ws.Range("C77").Value = 'Laura Siegemund'
$ws.Range("D77").Value = 'Hanyu Guo'
$ws.Range("E77").Value = 'Gana Hanyu Guo'
$ws.Range("F77").Value = 4
$ws.Range("G77").Value = "'"
$ws.Range("G77").ClearFormats()
$ws.Range("H77").Value = "'"
$ws.Range("H77").ClearFormats()

# Fila 78 - event_id 14359047
$ws.Range("A78").NumberFormat = "@"
$ws.Range("A78").Value = '14359047'
$ws.Range("A78").ClearFormats()
$ws.Range("B78").NumberFormat = "@"
$ws.Range("B78").Value = '2025-08-06'
$ws.Range("B78").ClearFormats()
$ws.Range("C78").Value = 'Iva Jovic'
$ws.Range("D78").Value = 'Varvara Gracheva'
$ws.Range("E78").Value = 'Gana Varvara Gracheva'
$ws.Range("F78").Value = 2.5
$ws.Range("G78").Value = "'"
$ws.Range("G78").ClearFormats()
$ws.Range("H78").Value = "'"
$ws.Range("H78").ClearFormats()

# Fila 79 - event_id 14359060
$ws.Range("A79").NumberFormat = "@"
$ws.Range("A79").Value = '14359060'
$ws.Range("A79").ClearFormats()
$ws.Range("B79").NumberFormat = "@"
$ws.Range("B79").Value = '2025-08-06'
$ws.Range("B79").ClearFormats()
$ws.Range("C79").Value = 'Dalma Galfi'
$ws.Range("D79").Value = 'Aoi Ito'
$ws.Range("E79").Value = 'Gana Aoi Ito'
$ws.Range("F79").Value = 2.1
$ws.Range("G79").Value = "'"
$ws.Range("G79").ClearFormats()
$ws.Range("H79").Value = "'"
$ws.Range("H79").ClearFormats()

# Fila 80 - event_id 14311067
$ws.Range("A80").NumberFormat = "@"
$ws.Range("A80").Value = '14311067'
$ws.Range("A80").ClearFormats()
$ws.Range("B80").NumberFormat = "@"
$ws.Range("B80").Value = '2025-08-07'
$ws.Range("B80").ClearFormats()
$ws.Range("C80").Value = 'Carlos Taberner'
$ws.Range("D80").Value = 'Max Alcala Gurri'
$ws.Range("E80").Value = 'Gana Max Alcala Gurri'
$ws.Range("F80").Value = 4.33
$ws.Range("G80").Value = "'"
$ws.Range("G80").ClearFormats()
$ws.Range("H80").Value = "'"
$ws.Range("H80").ClearFormats()

# Fila 81 - event_id 14310240
$ws.Range("A81").NumberFormat = "@"
$ws.Range("A81").Value = '14310240'
$ws.Range("A81").ClearFormats()
$ws.Range("B81").NumberFormat = "@"
$ws.Range("B81").Value = '2025-08-07'
$ws.Range("B81").ClearFormats()
$ws.Range("C81").Value = 'Jerome Kym'
$ws.Range("D81").Value = 'Raphael Collignon'
$ws.Range("E81").Value = 'Gana Jerome Kym'
$ws.Range("F81").Value = 2.25
$ws.Range("G81").Value = "'"
$ws.Range("G81").ClearFormats()
$ws.Range("H81").Value = "'"
$ws.Range("H81").ClearFormats()

# Fila 82 - event_id 14310237
$ws.Range("A82").NumberFormat = "@"
$ws.Range("A82").Value = '14310237'
$ws.Range("A82").ClearFormats()
$ws.Range("B82").NumberFormat = "@"
$ws.Range("B82").Value = '2025-08-07'
$ws.Range("B82").ClearFormats()
$ws.Range("C82").Value = 'Vilius Gaubas'
$ws.Range("D82").Value = 'Martin Krumich'
$ws.Range("E82").Value = 'Gana Martin Krumich'
$ws.Range("F82").Value = 4.33
$ws.Range("G82").Value = "'"
$ws.Range("G82").ClearFormats()
$ws.Range("H82").Value = "'"
$ws.Range("H82").ClearFormats()

# Fila 83 - event_id 14310263
$ws.Range("A83").NumberFormat = "@"
$ws.Range("A83").Value = '14310263'
$ws.Range("A83").ClearFormats()
$ws.Range("B83").NumberFormat = "@"
$ws.Range("B83").Value = '2025-08-06'
$ws.Range("B83").ClearFormats()
$ws.Range("C83").Value = 'Alex Bolt'
$ws.Range("D83").Value = 'Garrett Johns'
$ws.Range("E83").Value = 'Gana Garrett Johns'
$ws.Range("F83").Value = 3.4
$ws.Range("G83").Value = "'"
$ws.Range("G83").ClearFormats()
$ws.Range("H83").Value = "'"
$ws.Range("H83").ClearFormats()

# Fila 84 - event_id 14310257
$ws.Range("A84").NumberFormat = "@"
$ws.Range("A84").Value = '14310257'
$ws.Range("A84").ClearFormats()
$ws.Range("B84").NumberFormat = "@"
$ws.Range("B84").Value = '2025-08-06'
$ws.Range("B84").ClearFormats()
$ws.Range("C84").Value = 'Yibing Wu'
$ws.Range("D84").Value = 'Naoki Nakagawa'
$ws.Range("E84").Value = 'Gana Naoki Nakagawa'
$ws.Range("F84").Value = 10
$ws.Range("G84").Value = "'"
$ws.Range("G84").ClearFormats()
$ws.Range("H84").Value = "'"
$ws.Range("H84").ClearFormats()

